$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column keeps its original Text storage type (the source data
# stores prices/percentages as text, not numbers) so assigning numeric-looking
# strings like "1.00" or "686.11" does not get reinterpreted as a Number.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.580.15"
$ws.Range("E2").Value = "  +2.14%  "

$ws.Range("D3").Value = "3.816.41"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "686.11"
$ws.Range("E5").Value = "  +9.32%  "

$ws.Range("D6").Value = "170.06"
$ws.Range("E6").Value = "  +2.83%  "

$ws.Range("D7").Value = "3.813.99"
$ws.Range("E7").Value = "  +1.03%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("E10").Value = "  +1.63%  "

$ws.Range("D11").Value = "7.20"
$ws.Range("E11").Value = "  +6.30%  "

$ws.Range("E12").Value = "  +0.71%  "

$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("D14").Value = "35.85"
$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("D15").Value = "4.461.46"
$ws.Range("E15").Value = "  +1.24%  "

$ws.Range("D16").Value = "3.817.12"
$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("D17").Value = "70.673.84"
$ws.Range("E17").Value = "  +2.25%  "

$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("E20").Value = "  +0.51%  "

$ws.Range("D21").Value = "11.24"
$ws.Range("E21").Value = "  +17.71%  "

$ws.Range("D22").Value = "477.57"
$ws.Range("E22").Value = "  +2.19%  "

$ws.Range("E23").Value = "  +1.16%  "

$ws.Range("D24").Value = "83.30"
$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("E25").Value = "  -1.57%  "

$ws.Range("D26").Value = "12.29"
$ws.Range("E26").Value = "  +2.36%  "

$ws.Range("D27").Value = "10.34"
$ws.Range("E27").Value = "  +3.16%  "

$ws.Range("E28").Value = "  -1.96%  "

$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("D30").Value = "3.968.15"
$ws.Range("E30").Value = "  +1.17%  "

$ws.Range("D31").Value = "2.94"
$ws.Range("E31").Value = "  +10.23%  "

$ws.Range("E32").Value = "  +2.84%  "

$ws.Range("D33").Value = "7.41"
$ws.Range("E33").Value = "  +3.84%  "

$ws.Range("D34").Value = "29.58"
$ws.Range("E34").Value = "  +2.87%  "

$ws.Range("E35").Value = "  +3.64%  "

$ws.Range("D36").Value = "9.15"
$ws.Range("E36").Value = "  +2.32%  "

$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").Value = "3.773.12"
$ws.Range("E38").Value = "  +1.31%  "

$ws.Range("E39").Value = "  +1.32%  "

$ws.Range("D40").Value = "3.39"
$ws.Range("E40").Value = "  +2.12%  "

$ws.Range("E41").Value = "  +2.08%  "

$ws.Range("D42").Value = "0.964"
$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("E44").Value = "  +11.62%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "46.00"

$ws.Range("D47").Value = "159.52"
$ws.Range("E47").Value = "  +3.47%  "

$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "0.000299"
$ws.Range("E48").Value = "  +10.43%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "48.18"
$ws.Range("E49").Value = "  +3.04%  "

$ws.Range("E50").Value = "  +6.27%  "

$ws.Range("E51").Value = "  +1.62%  "
